# edit.ps1 - apply the text edits described by the target diff
$p = $ppt.ActivePresentation

# --- Slide 1: title text change (second run of the single title shape) ---
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$para1 = $sh1.TextFrame.TextRange.Paragraphs(1)
$para1.Runs(2).Text = "BB22:The Most Advanced QKD"

# --- Slide 13: "If Bob Measures with Z Basis" bullet list, 2nd bullet ---
# Simple in-place text fix: |1> -> |0> (no run split in the diff)
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(2)
$tr13 = $sh13.TextFrame.TextRange
$run13 = $tr13.Paragraphs(2).Runs(1)
$run13.Text = "If he measures “0”s vast majority of the time, he is receiving the state of |0⟩"

# --- Slide 14: "If Bob Measures with X-Basis" bullet list ---
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(2)
$tr14 = $sh14.TextFrame.TextRange

# Paragraph 1: split into two runs (text content unchanged)
$run14_1a = $tr14.Paragraphs(1).Runs(1)
$run14_1a.Text = "If he measures “0”s vast majority of the time, he "
$run14_1b = $run14_1a.InsertAfter("is receiving |+⟩")

# Paragraph 2: text changes from "0" to "1" AND splits into two runs
$run14_2a = $tr14.Paragraphs(2).Runs(1)
$run14_2a.Text = "If he measures “1”s vast majority of the time, he "
$run14_2b = $run14_2a.InsertAfter("is receiving |-⟩")

# Paragraph 3: splits into three runs (text content unchanged)
$run14_3a = $tr14.Paragraphs(3).Runs(1)
$run14_3a.Text = "If he measures “0”s and “1”s 50:50 distributed, "
$run14_3b = $run14_3a.InsertAfter("he is receiving the state of |0⟩ or |1⟩, depending ")
$run14_3c = $run14_3b.InsertAfter("on phase which cannot be measured directly")

# --- Slide 3: "Benefits" bullet list, 2nd bullet splits into two runs ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange
$run3a = $tr3.Paragraphs(2).Runs(1)
$run3a.Text = "Presence of eavesdropping can be detected "
$run3b = $run3a.InsertAfter("apparently")
